$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.922420666666667
$ws.Range("H2").Value = 17.767262
$ws.Range("I2").Value = 0.5833698282960434
$ws.Range("J2").Value = 0.6311054116979437
$ws.Range("M2").Value = 30.46625333333334
$ws.Range("N2").Value = 91.39876000000001
$ws.Range("O2").Value = 0.2185380492512374
$ws.Range("P2").Value = 0.2331534018544084
$ws.Range("Q2").Value = 180.4339683772356
$ws.Range("R2").Value = 1623.90571539512
$ws.Range("S2").Value = 0.1274885042678467
$ws.Range("T2").Value = 0.1471443736661025
$ws.Range("G3").Value = 5.922420666666667
$ws.Range("H3").Value = 17.767262
$ws.Range("I3").Value = 0.5833698282960434
$ws.Range("J3").Value = 0.6311054116979437
$ws.Range("O3").Value = 0.2491807703757967
$ws.Range("P3").Value = 0.2658454419670822
$ws.Range("Q3").Value = 205.7338545678773
$ws.Range("R3").Value = 1851.604691110896
$ws.Range("S3").Value = 0.1453645432288043
$ws.Range("T3").Value = 0.1677764971006572
$ws.Range("G4").Value = 5.922420666666667
$ws.Range("H4").Value = 17.767262
$ws.Range("I4").Value = 0.5833698282960434
$ws.Range("J4").Value = 0.6311054116979437
$ws.Range("M4").Value = 23.69037333333334
$ws.Range("N4").Value = 71.07112000000001
$ws.Range("O4").Value = 0.1699338582153697
$ws.Range("P4").Value = 0.181298667526812
$ws.Range("Q4").Value = 140.3043566303822
$ws.Range("R4").Value = 1262.73920967344
$ws.Range("S4").Value = 0.09913428568878443
$ws.Range("T4").Value = 0.1144185702097973
$ws.Range("G5").Value = 5.922420666666667
$ws.Range("H5").Value = 17.767262
$ws.Range("I5").Value = 0.5833698282960434
$ws.Range("J5").Value = 0.6311054116979437
$ws.Range("M5").Value = 26.2168665
$ws.Range("N5").Value = 52.433733
$ws.Range("O5").Value = 0.18805669340777
$ws.Range("P5").Value = 0.1337556791894743
$ws.Range("Q5").Value = 155.267311974841
$ws.Range("R5").Value = 931.603871849046
$ws.Range("S5").Value = 0.1097066009432125
$ws.Range("T5").Value = 0.08441393298181125
$ws.Range("G6").Value = 5.922420666666667
$ws.Range("H6").Value = 17.767262
$ws.Range("I6").Value = 0.5833698282960434
$ws.Range("J6").Value = 0.6311054116979437
$ws.Range("M6").Value = 24.297748
$ws.Range("N6").Value = 72.893244
$ws.Range("O6").Value = 0.1742906287498262
$ws.Range("P6").Value = 0.1859468094622229
$ws.Range("Q6").Value = 143.9014849086587
$ws.Range("R6").Value = 1295.113364177928
$ws.Range("S6").Value = 0.1016758941673956
$ws.Range("T6").Value = 0.1173520377395753
$ws.Range("I7").Value = 0.1897160182974547
$ws.Range("J7").Value = 0.2052399696141807
$ws.Range("M7").Value = 30.46625333333334
$ws.Range("N7").Value = 91.39876000000001
$ws.Range("O7").Value = 0.2185380492512374
$ws.Range("P7").Value = 0.2331534018544084
$ws.Range("Q7").Value = 58.67841013671113
$ws.Range("R7").Value = 528.1056912304001
$ws.Range("S7").Value = 0.04146016855043781
$ws.Range("T7").Value = 0.04785239711204165
$ws.Range("I8").Value = 0.1897160182974547
$ws.Range("J8").Value = 0.2052399696141807
$ws.Range("O8").Value = 0.2491807703757967
$ws.Range("P8").Value = 0.2658454419670822
$ws.Range("Q8").Value = 66.90611311114667
$ws.Range("S8").Value = 0.04727358359198849
$ws.Range("T8").Value = 0.05456211043139239
$ws.Range("I9").Value = 0.1897160182974547
$ws.Range("J9").Value = 0.2052399696141807
$ws.Range("M9").Value = 23.69037333333334
$ws.Range("N9").Value = 71.07112000000001
$ws.Range("O9").Value = 0.1699338582153697
$ws.Range("P9").Value = 0.181298667526812
$ws.Range("Q9").Value = 45.62797491164446
$ws.Range("R9").Value = 410.6517742048001
$ws.Range("S9").Value = 0.03223917495454415
$ws.Range("T9").Value = 0.03720973301429435
$ws.Range("I10").Value = 0.1897160182974547
$ws.Range("J10").Value = 0.2052399696141807
$ws.Range("M10").Value = 26.2168665
$ws.Range("N10").Value = 52.433733
$ws.Range("O10").Value = 0.18805669340777
$ws.Range("P10").Value = 0.1337556791894743
$ws.Range("Q10").Value = 50.49403443722001
$ws.Range("R10").Value = 302.9642066233201
$ws.Range("S10").Value = 0.03567736708750732
$ws.Range("T10").Value = 0.02745201153257181
$ws.Range("I11").Value = 0.1897160182974547
$ws.Range("J11").Value = 0.2052399696141807
$ws.Range("M11").Value = 24.297748
$ws.Range("N11").Value = 72.893244
$ws.Range("O11").Value = 0.1742906287498262
$ws.Range("P11").Value = 0.1859468094622229
$ws.Range("Q11").Value = 46.79778661797334
$ws.Range("R11").Value = 421.18007956176
$ws.Range("S11").Value = 0.03306572411297691
$ws.Range("T11").Value = 0.03816371752388049
$ws.Range("G12").Value = 2.303652
$ws.Range("H12").Value = 4.607303999999999
$ws.Range("I12").Value = 0.2269141534065018
$ws.Range("J12").Value = 0.1636546186878756
$ws.Range("M12").Value = 30.46625333333334
$ws.Range("N12").Value = 91.39876000000001
$ws.Range("O12").Value = 0.2185380492512374
$ws.Range("P12").Value = 0.2331534018544084
$ws.Range("Q12").Value = 70.18364542383999
$ws.Range("R12").Value = 421.10187254304
$ws.Range("S12").Value = 0.04958937643295294
$ws.Range("T12").Value = 0.03815663107626425
$ws.Range("G13").Value = 2.303652
$ws.Range("H13").Value = 4.607303999999999
$ws.Range("I13").Value = 0.2269141534065018
$ws.Range("J13").Value = 0.1636546186878756
$ws.Range("O13").Value = 0.2491807703757967
$ws.Range("P13").Value = 0.2658454419670822
$ws.Range("Q13").Value = 80.02457647267198
$ws.Range("R13").Value = 480.1474588360319
$ws.Range("S13").Value = 0.05654264355500382
$ws.Range("T13").Value = 0.04350683443503261
$ws.Range("G14").Value = 2.303652
$ws.Range("H14").Value = 4.607303999999999
$ws.Range("I14").Value = 0.2269141534065018
$ws.Range("J14").Value = 0.1636546186878756
$ws.Range("M14").Value = 23.69037333333334
$ws.Range("N14").Value = 71.07112000000001
$ws.Range("O14").Value = 0.1699338582153697
$ws.Range("P14").Value = 0.181298667526812
$ws.Range("Q14").Value = 54.57437591008
$ws.Range("R14").Value = 327.44625546048
$ws.Range("S14").Value = 0.03856039757204114
$ws.Range("T14").Value = 0.02967036430272036
$ws.Range("G15").Value = 2.303652
$ws.Range("H15").Value = 4.607303999999999
$ws.Range("I15").Value = 0.2269141534065018
$ws.Range("J15").Value = 0.1636546186878756
$ws.Range("M15").Value = 26.2168665
$ws.Range("N15").Value = 52.433733
$ws.Range("O15").Value = 0.18805669340777
$ws.Range("P15").Value = 0.1337556791894743
$ws.Range("Q15").Value = 60.394536946458
$ws.Range("R15").Value = 241.578147785832
$ws.Range("S15").Value = 0.0426727253770502
$ws.Range("T15").Value = 0.02188973467509124
$ws.Range("G16").Value = 2.303652
$ws.Range("H16").Value = 4.607303999999999
$ws.Range("I16").Value = 0.2269141534065018
$ws.Range("J16").Value = 0.1636546186878756
$ws.Range("M16").Value = 24.297748
$ws.Range("N16").Value = 72.893244
$ws.Range("O16").Value = 0.1742906287498262
$ws.Range("P16").Value = 0.1859468094622229
$ws.Range("Q16").Value = 55.97355577569598
$ws.Range("R16").Value = 335.8413346541759
$ws.Range("S16").Value = 0.03954901046945372
$ws.Range("T16").Value = 0.03043105419876716

$wb.Save()
